$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1596
$ws1.Range("F4").Value = 5215
$ws1.Range("F5").Value = 560
$ws1.Range("F6").Value = 10367
$ws1.Range("F8").Value = 565
$ws1.Range("F9").Value = 117
$ws1.Range("F10").Value = 121
$ws1.Range("F11").Value = 823
$ws1.Range("F12").Value = 85

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1596
$ws4.Range("F6").Value = 5215
$ws4.Range("F7").Value = 560
$ws4.Range("F9").Value = 10367
$ws4.Range("F11").Value = 565
$ws4.Range("F12").Value = 117
$ws4.Range("F15").Value = 121
$ws4.Range("F16").Value = 823
$ws4.Range("F18").Value = 85
